# Últimas simulações para o parâmetro de impacto
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: object / simulation parameters ---
$ws.Range("B2").Value = "WASP74b_Hellier2015_TL+"
$ws.Range("C2").Value = 700
$ws.Range("E2").Value = 1600
$ws.Range("F2").Value = 1.59
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("L2").Formula = "=SQRT(0.08)"
$ws.Range("N2").Value = 0.038
$ws.Range("O2").Value = 1.48
$ws.Range("P2").Value = 1.56
$ws.Range("Q2").Value = 0.95
$ws.Range("R2").Value = 2.13775

# anguloInclinacao (S2): value change + simplified number format (12 decimals -> 3 decimals)
$ws.Range("S2").NumberFormat = "0.000"
$ws.Range("S2").Value = 79.57

$ws.Range("W2").Value = 1

# --- Header label change: semiEixoRaioStar -> semiEixoUA (column N) ---
$ws.Range("N1").Value = "semiEixoUA"

# --- Row 3: clear the second mancha (lat/longt/r) ---
$ws.Range("J3:L3").ClearContents()

# --- Selection moves to L2 ---
$ws.Range("L2").Select()
